# edit.ps1 - applies the week1/演示.pptx OOXML diff via PowerPoint COM-interop.
#
# Summary of the three changes:
#  1. notesSlide2 (Slide 2 notes body): delete the 2nd paragraph
#     ("这里的chip方法把基因切成小的片段...等位基因也是一个非常小的区域。").
#  2. notesSlide3 (Slide 3 notes body): replace the first two paragraphs
#     ("act in the expected direction????" and "作者用杂合子样本，那要是
#     没有等位基因怎么办？方法不具有普适性。") with one new paragraph about
#     heterozygous SNPs.
#  3. Slide 2, subtitle shape: shorten "...的等位基因（SNP）会在..." to
#     "...的等位SNP会在..." (drop "基因（" and the "）" run).
#
# Notes: this COM host only supports whole-shape TextRange.Text assignment
# for NotesPage placeholders (sub-range edits / Font changes raise "could not
# set font properties" / "could not apply the edit" there), so items 1 & 2
# are done by reconstructing the full notes body text. Item 3 is a normal
# slide shape, where precise Characters(start,len) sub-range edits work and
# keep every other run's rPr untouched, so that one is done surgically.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 2 notes: drop the second paragraph entirely.
# ---------------------------------------------------------------------
$notes2 = $p.Slides.Item(2).NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes2.Text = '测的都是提前知道与TF-binding有关的样本。'

# ---------------------------------------------------------------------
# 2) Slide 3 notes: merge/replace paragraphs 1-2 into a single new one,
#    keep paragraphs 3-6 unchanged.
# ---------------------------------------------------------------------
$notes3 = $p.Slides.Item(3).NotesPage.Shapes.Item(2).TextFrame.TextRange

$newPara1 = '作者用杂合子样本，杂合指的是SNP的杂合，等位SNP指的不是处于同源染色体上，而是SNP的ATCG的碱基差别。'

$keepParas = @(
    '混合起来测就可以提升效率的话，此前的人为什么不做？',
    '问题：这里的frequency肯定不能是相对于全部基因的，否则结果会由少数大量且能与转录因子结合的基因主导。',
    '          但是，如果只计算相对于等位基因的比例，当一对等位基因都与TF有很强的结合性时，反而前后不会出现很大的frequency差异。导致无法筛选出来。',
    '          此外，这种混合方法建立在比较frequency上，转而比较基因前后量也是不合理的，因为pool中有其他基因存在，而且基因的量也不固定。少量基因A，在有一些其他基因的pool里剩下20%，和大量基因B，在有另一些其他基因的pool里剩下20%是不可比较的。'
)

$lf = [char]10
$newNotes3Text = $newPara1
foreach ($para in $keepParas) {
    $newNotes3Text = $newNotes3Text + $lf + $para
}
$notes3.Text = $newNotes3Text

# ---------------------------------------------------------------------
# 3) Slide 2, subtitle shape ("副标题 6"): trim "等位基因（" -> "等位" and
#    remove the standalone "）" run, without touching surrounding runs.
# ---------------------------------------------------------------------
$sub = $p.Slides.Item(2).Shapes.Item(2).TextFrame.TextRange

$oldPhrase = '在顺式作用时促进转录因子结合的等位基因（'
$newPhrase = '在顺式作用时促进转录因子结合的等位'

# Locate the run precisely instead of hard-coding offsets: scan the whole
# text for the old phrase.
$fullText = $sub.Text
$idx = $fullText.IndexOf($oldPhrase)
if ($idx -lt 0) {
    throw "could not find target phrase in slide 2 subtitle shape"
}
$startPos = $idx + 1   # TextRange is 1-based
$run = $sub.Characters($startPos, $oldPhrase.Length)
$run.Text = $newPhrase

# After trimming, the following "）" run shifts left by the number of
# characters removed (3, matching the ideographic parenthesis already
# accounted for: "基因（" minus "" is 3 chars removed). Find it again
# relative to "SNP" so we don't depend on a fixed offset.
$fullText2 = $sub.Text
$afterNewPhrase = $idx + $newPhrase.Length  # 0-based index right after new phrase
$snpIdx = $fullText2.IndexOf('SNP', $afterNewPhrase)
if ($snpIdx -lt 0) {
    throw "could not find SNP run after trimmed phrase"
}
$parenZeroIdx = $snpIdx + 3   # 0-based index of the char right after "SNP"
$parenPos = $parenZeroIdx + 1  # 1-based
$parenRun = $sub.Characters($parenPos, 1)
if ($parenRun.Text -eq '）') {
    $parenRun.Text = ''
}
